$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fn1"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 35.73885133333334
$ws.Cells.Item(2, 8).Value = 107.216554
$ws.Cells.Item(2, 9).Value = 0.01949729408921566
$ws.Cells.Item(2, 10).Value = 0.01949729408921566
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.06861733333333334
$ws.Cells.Item(2, 14).Value = 0.205852
$ws.Cells.Item(2, 15).Value = 0.01654048691795588
$ws.Cells.Item(2, 16).Value = 0.01654048691795588
$ws.Cells.Item(2, 17).Value = 2.452304674889778
$ws.Cells.Item(2, 18).Value = 22.070742074008
$ws.Cells.Item(2, 19).Value = 0.0003224947378182102
$ws.Cells.Item(2, 20).Value = 0.0003224947378182102

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fn1"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 35.73885133333334
$ws.Cells.Item(3, 8).Value = 107.216554
$ws.Cells.Item(3, 9).Value = 0.01949729408921566
$ws.Cells.Item(3, 10).Value = 0.01949729408921566
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.776574666666666
$ws.Cells.Item(3, 14).Value = 11.329724
$ws.Cells.Item(3, 15).Value = 0.9103586635352137
$ws.Cells.Item(3, 16).Value = 0.9103586635352137
$ws.Cells.Item(3, 17).Value = 134.9704405612329
$ws.Cells.Item(3, 18).Value = 1214.733965051096
$ws.Cells.Item(3, 19).Value = 0.01774953058961139
$ws.Cells.Item(3, 20).Value = 0.01774953058961139

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fn1"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 35.73885133333334
$ws.Cells.Item(4, 8).Value = 107.216554
$ws.Cells.Item(4, 9).Value = 0.01949729408921566
$ws.Cells.Item(4, 10).Value = 0.01949729408921566
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.303255
$ws.Cells.Item(4, 14).Value = 0.909765
$ws.Cells.Item(4, 15).Value = 0.07310084954683041
$ws.Cells.Item(4, 16).Value = 0.07310084954683042
$ws.Cells.Item(4, 17).Value = 10.83798536109
$ws.Cells.Item(4, 18).Value = 97.54186824981001
$ws.Cells.Item(4, 19).Value = 0.00142526876178606
$ws.Cells.Item(4, 20).Value = 0.00142526876178606

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fn1"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1689.289306666667
$ws.Cells.Item(5, 8).Value = 5067.86792
$ws.Cells.Item(5, 9).Value = 0.9215900675332435
$ws.Cells.Item(5, 10).Value = 0.9215900675332435
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.06861733333333334
$ws.Cells.Item(5, 14).Value = 0.205852
$ws.Cells.Item(5, 15).Value = 0.01654048691795588
$ws.Cells.Item(5, 16).Value = 0.01654048691795588
$ws.Cells.Item(5, 17).Value = 115.9145274519822
$ws.Cells.Item(5, 18).Value = 1043.23074706784
$ws.Cells.Item(5, 19).Value = 0.01524354845575169
$ws.Cells.Item(5, 20).Value = 0.01524354845575169

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fn1"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1689.289306666667
$ws.Cells.Item(6, 8).Value = 5067.86792
$ws.Cells.Item(6, 9).Value = 0.9215900675332435
$ws.Cells.Item(6, 10).Value = 0.9215900675332435
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.776574666666666
$ws.Cells.Item(6, 14).Value = 11.329724
$ws.Cells.Item(6, 15).Value = 0.9103586635352137
$ws.Cells.Item(6, 16).Value = 0.9103586635352137
$ws.Cells.Item(6, 17).Value = 6379.727200228231
$ws.Cells.Item(6, 18).Value = 57417.54480205407
$ws.Cells.Item(6, 19).Value = 0.8389775022068908
$ws.Cells.Item(6, 20).Value = 0.8389775022068908

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fn1"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1689.289306666667
$ws.Cells.Item(7, 8).Value = 5067.86792
$ws.Cells.Item(7, 9).Value = 0.9215900675332435
$ws.Cells.Item(7, 10).Value = 0.9215900675332435
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.303255
$ws.Cells.Item(7, 14).Value = 0.909765
$ws.Cells.Item(7, 15).Value = 0.07310084954683041
$ws.Cells.Item(7, 16).Value = 0.07310084954683042
$ws.Cells.Item(7, 17).Value = 512.2854286931999
$ws.Cells.Item(7, 18).Value = 4610.5688582388
$ws.Cells.Item(7, 19).Value = 0.0673690168706009
$ws.Cells.Item(7, 20).Value = 0.06736901687060091

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Fn1"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 93.641553
$ws.Cells.Item(8, 8).Value = 280.924659
$ws.Cells.Item(8, 9).Value = 0.05108605424341119
$ws.Cells.Item(8, 10).Value = 0.05108605424341119
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06861733333333334
$ws.Cells.Item(8, 14).Value = 0.205852
$ws.Cells.Item(8, 15).Value = 0.01654048691795588
$ws.Cells.Item(8, 16).Value = 0.01654048691795588
$ws.Cells.Item(8, 17).Value = 6.425433656052
$ws.Cells.Item(8, 18).Value = 57.82890290446801
$ws.Cells.Item(8, 19).Value = 0.0008449882119031275
$ws.Cells.Item(8, 20).Value = 0.0008449882119031275

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Fn1"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 93.641553
$ws.Cells.Item(9, 8).Value = 280.924659
$ws.Cells.Item(9, 9).Value = 0.05108605424341119
$ws.Cells.Item(9, 10).Value = 0.05108605424341119
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.776574666666666
$ws.Cells.Item(9, 14).Value = 11.329724
$ws.Cells.Item(9, 15).Value = 0.9103586635352137
$ws.Cells.Item(9, 16).Value = 0.9103586635352137
$ws.Cells.Item(9, 17).Value = 353.644316807124
$ws.Cells.Item(9, 18).Value = 3182.798851264116
$ws.Cells.Item(9, 19).Value = 0.04650663206631924
$ws.Cells.Item(9, 20).Value = 0.04650663206631924

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Fn1"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 93.641553
$ws.Cells.Item(10, 8).Value = 280.924659
$ws.Cells.Item(10, 9).Value = 0.05108605424341119
$ws.Cells.Item(10, 10).Value = 0.05108605424341119
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.303255
$ws.Cells.Item(10, 14).Value = 0.909765
$ws.Cells.Item(10, 15).Value = 0.07310084954683041
$ws.Cells.Item(10, 16).Value = 0.07310084954683042
$ws.Cells.Item(10, 17).Value = 28.397269155015
$ws.Cells.Item(10, 18).Value = 255.575422395135
$ws.Cells.Item(10, 19).Value = 0.003734433965188818
$ws.Cells.Item(10, 20).Value = 0.003734433965188819

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Fn1"
$ws.Cells.Item(11, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 14.34625366666667
$ws.Cells.Item(11, 8).Value = 43.038761
$ws.Cells.Item(11, 9).Value = 0.007826584134129748
$ws.Cells.Item(11, 10).Value = 0.007826584134129748
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.06861733333333334
$ws.Cells.Item(11, 14).Value = 0.205852
$ws.Cells.Item(11, 15).Value = 0.01654048691795588
$ws.Cells.Item(11, 16).Value = 0.01654048691795588
$ws.Cells.Item(11, 17).Value = 0.9844016699302223
$ws.Cells.Item(11, 18).Value = 8.859615029372
$ws.Cells.Item(11, 19).Value = 0.0001294555124828542
$ws.Cells.Item(11, 20).Value = 0.0001294555124828542

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Fn1"
$ws.Cells.Item(12, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 14.34625366666667
$ws.Cells.Item(12, 8).Value = 43.038761
$ws.Cells.Item(12, 9).Value = 0.007826584134129748
$ws.Cells.Item(12, 10).Value = 0.007826584134129748
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.776574666666666
$ws.Cells.Item(12, 14).Value = 11.329724
$ws.Cells.Item(12, 15).Value = 0.9103586635352137
$ws.Cells.Item(12, 16).Value = 0.9103586635352137
$ws.Cells.Item(12, 17).Value = 54.17969815910711
$ws.Cells.Item(12, 18).Value = 487.617283431964
$ws.Cells.Item(12, 19).Value = 0.007124998672392265
$ws.Cells.Item(12, 20).Value = 0.007124998672392265

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Fn1"
$ws.Cells.Item(13, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 14.34625366666667
$ws.Cells.Item(13, 8).Value = 43.038761
$ws.Cells.Item(13, 9).Value = 0.007826584134129748
$ws.Cells.Item(13, 10).Value = 0.007826584134129748
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.303255
$ws.Cells.Item(13, 14).Value = 0.909765
$ws.Cells.Item(13, 15).Value = 0.07310084954683041
$ws.Cells.Item(13, 16).Value = 0.07310084954683042
$ws.Cells.Item(13, 17).Value = 4.350573155685
$ws.Cells.Item(13, 18).Value = 39.155158401165
$ws.Cells.Item(13, 19).Value = 0.0005721299492546287
$ws.Cells.Item(13, 20).Value = 0.0005721299492546288
